# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets
# to reflect the latest scrape, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 817
    3  = 6
    4  = 1140
    5  = 45
    6  = 12277
    9  = 490
    10 = 435
    12 = 893
    13 = 13591
    14 = 13749
    19 = 1026
    20 = 102
    22 = 3622
    23 = 208
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
